$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to retain text formatting so numeric-looking
# strings (e.g. "1.003") are not auto-converted to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.105.73"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.656.40"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "218.89"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "0.5287"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").Value = "0.2608"
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("D9").Value = "0.06343"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "20.39"
$ws.Range("E10").Value = "  -2.71%  "
$ws.Range("D11").Value = "0.07796"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "4.498"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "1.665.91"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "0.5472"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "0.0₅8169"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "65.20"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "26.127.97"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "4.553"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("D20").Value = "192.65"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "10.08"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").Value = "6.036"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "140.61"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "0.1240"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "7.262"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "1.433"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Value = "0.05907"
$ws.Range("E29").Value = "  -3.76%  "
$ws.Range("D30").Value = "1.282"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").Value = "3.523"
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("D32").Value = "3.238"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").Value = "1.558"
$ws.Range("E33").Value = "  -5.05%  "
$ws.Range("D34").Value = "0.9475"
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("D35").Value = "2.412"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").Value = "2.774"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "0.5629"
$ws.Range("E37").Value = "  -2.40%  "
$ws.Range("D38").Value = "0.01611"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").Value = "5.790"
$ws.Range("E39").Value = "  -4.50%  "
$ws.Range("D40").Value = "0.8449"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "101.59"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("D43").Value = "1.008.97"
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").Value = "1.799.10"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").Value = "57.01"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈106"
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "1.005"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").Value = "0.4289"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").Value = "1.480"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "0.05157"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "7.748"
$ws.Range("E51").Value = "  -4.64%  "

# Reset style on the price column so only the value/type changed (no
# leftover quote-prefix/text-format styling), matching the source data
# which stores these as plain inline strings with the default style.
$priceRange.Style = "Normal"

